$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content changes in column G ---
# Remove "checkGameOver" (row 5) entirely.
$ws.Range("G5").ClearContents()

# "getGameSquares" moves from G23 up to G20.
$ws.Range("G20").Value = $ws.Range("G23").Value()
$ws.Range("G23").ClearContents()

# --- Formatting changes ---
# Functions considered "done"/obsolete are now shown struck-through.
$ws.Range("G4").Font.Strikethrough = $true
$ws.Range("G6").Font.Strikethrough = $true
$ws.Range("G7").Font.Strikethrough = $true
$ws.Range("G8").Font.Strikethrough = $true
$ws.Range("G17").Font.Strikethrough = $true
$ws.Range("G18").Font.Strikethrough = $true

# G13 loses its special (Arial Unicode MS) font, back to the default style.
$ws.Range("G13").Style = "Normal"

# Row 13 no longer has the taller custom row height (back to default).
$ws.Rows.Item(13).AutoFit()

# New column width for column G (COM ColumnWidth units are offset by 5/6
# from the raw OOXML column width, hence the + 1/12 vs. the desired 25).
$ws.Columns.Item(7).ColumnWidth = 25 + 1/12

# Update the active selection shown when the workbook is reopened.
$ws.Range("H9").Select()
